$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "26.389.58"
$ws.Range("E2").Value2 = "  -0.88%  "
$ws.Range("D3").Value2 = "1.593.23"
$ws.Range("E3").Value2 = "  -0.22%  "
$ws.Range("E4").Value2 = "  -0.70%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "210.37"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  -0.27%  "
$ws.Range("E7").Value2 = "  -0.66%  "
$ws.Range("E8").Value2 = "  -0.67%  "
$ws.Range("E9").Value2 = "  -0.11%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "19.55"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  +0.15%  "
$ws.Range("E11").Value2 = "  -0.03%  "
$ws.Range("D12").Value2 = "1.818.58"
$ws.Range("E12").Value2 = "  -0.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = "4.07"
$c.Style = "Normal"
$ws.Range("E13").Value2 = "  +0.81%  "
$ws.Range("D14").Value2 = "1.561.51"
$ws.Range("E14").Value2 = "  -2.31%  "
$ws.Range("E15").Value2 = "  -0.55%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = "64.48"
$c.Style = "Normal"
$ws.Range("E16").Value2 = "  -0.14%  "
$ws.Range("D17").Value2 = "26.391.38"
$ws.Range("E17").Value2 = "  -0.82%  "
$ws.Range("E18").Value2 = "  -1.19%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "7.49"
$c.Style = "Normal"
$ws.Range("E19").Value2 = "  +5.56%  "
$ws.Range("E20").Value2 = "  +1.46%  "
$ws.Range("E21").Value2 = "  -0.65%  "
$ws.Range("E22").Value2 = "  +0.10%  "
$ws.Range("B23").Value2 = "Avalanche"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = "8.94"
$c.Style = "Normal"
$ws.Range("E23").Value2 = "  +0.17%  "
$ws.Range("B24").Value2 = "Toncoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "2.15"
$c.Style = "Normal"
$ws.Range("E24").Value2 = "  -4.12%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "145.03"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  +0.87%  "
$ws.Range("E26").Value2 = "  -0.68%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "7.07"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  -0.85%  "
$ws.Range("E28").Value2 = "  +0.01%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "15.28"
$c.Style = "Normal"
$ws.Range("E29").Value2 = "  +0.25%  "
$ws.Range("E30").Value2 = "  -0.03%  "
$ws.Range("E31").Value2 = "  -0.39%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = "3.22"
$c.Style = "Normal"
$ws.Range("E32").Value2 = "  -0.53%  "
$ws.Range("E33").Value2 = "  +1.37%  "
$ws.Range("D34").Value2 = "1.311.68"
$ws.Range("E34").Value2 = "  +2.86%  "
$ws.Range("E35").Value2 = "  +3.42%  "
$ws.Range("E36").Value2 = "  -1.85%  "
$ws.Range("E37").Value2 = "  -0.41%  "
$ws.Range("E38").Value2 = "  +0.49%  "
$ws.Range("E39").Value2 = "  -13.18%  "
$ws.Range("E40").Value2 = "  -0.83%  "
$ws.Range("E41").Value2 = "  -0.66%  "
$ws.Range("E42").Value2 = "  +4.16%  "
$ws.Range("B43").Value2 = "TrustWalletToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "0.766"
$c.Style = "Normal"
$ws.Range("E43").Value2 = "  -1.37%  "
$ws.Range("B44").Value2 = "MXToken"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "2.14"
$c.Style = "Normal"
$ws.Range("E44").Value2 = "  -1.20%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = "62.64"
$c.Style = "Normal"
$ws.Range("E45").Value2 = "  +0.20%  "
$ws.Range("D46").Value2 = "1.729.79"
$ws.Range("E46").Value2 = "  -0.19%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "87.97"
$c.Style = "Normal"
$ws.Range("E47").Value2 = "  -2.11%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "1.50"
$c.Style = "Normal"
$ws.Range("E48").Value2 = "  -3.98%  "
$ws.Range("E49").Value2 = "  -1.59%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "0.0983"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -3.94%  "
$ws.Range("E51").Value2 = "  -1.43%  "
